# Apply updated cryptocurrency price/volume data to the worksheet.
# A leading apostrophe forces Excel to store a numeric-looking price
# string (e.g. "211.32") as text, matching the original inline-string layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.354.35"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.592.02"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'211.32"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'0.506"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "'19.46"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.06"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.579.94"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "'64.72"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "26.357.93"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("D19").Value = "'7.49"
$ws.Range("E19").Value = "  +3.46%  "
$ws.Range("D20").Value = "'211.90"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").Value = "'4.30"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("E23").Value = "  -1.65%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").Value = "'144.81"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "'7.08"
$ws.Range("E27").Value = "  +0.77%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +0.95%  "
$ws.Range("D34").Value = "1.340.27"
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("E35").Value = "  -1.25%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.599"
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'1.48"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -0.24%  "
$ws.Range("D39").Value = "'1.07"
$ws.Range("E39").Value = "  -15.23%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  +4.77%  "
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "1.728.75"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'61.62"
$ws.Range("E46").Value = "  -1.07%  "
$ws.Range("D47").Value = "'88.01"
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = "  -0.32%  "
